$d = $word.ActiveDocument

# 1. Update the date/day heading.
$d.Content.Find.Execute("2025-02-21 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-02-22 Saturday", 2)

# 2. Update the practice-problem table. The table only has content in
#    physical rows 1, 5, 9, 13, 17 (the other rows are blank spacer rows).
#    Each of those rows has 5 columns; the cell count per row does not
#    change, only the individual problem text.
$t = $d.Tables.Item(1)

$rows = @{
    1  = @("42÷7=", "75÷8=", "10÷2=", "62÷2=", "39÷3=")
    5  = @("54÷5=", "62÷6=", "29÷8=", "77÷9=", "25÷3=")
    9  = @("82÷2=", "57÷2=", "69÷7=", "39÷9=", "52÷9=")
    13 = @("30÷3=", "58÷2=", "91÷9=", "92÷6=", "18÷4=")
    17 = @("25÷4=", "33÷4=", "23÷9=", "60÷9=", "39÷5=")
}

foreach ($rowIndex in $rows.Keys) {
    $values = $rows[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
